# Saldo_guide.xlsx update
# - Rename the sheet to reflect the refreshed export file name
# - Bump the "Dt. Referencia" (reference date) for every data row from
#   2024-05-09 (serial 45421) to 2024-05-10 (serial 45422)
# - Refresh per-account balances: "Vl. Projetado" (projected value) gets
#   folded into "Saldo Previsto" (expected balance) for rows whose totals
#   moved, and "Vl. Total" is recalculated accordingly
# - Move the active cell selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "IClientBalance-20240510-091832-"

# Bump the reference date column (G) for every data row, 2024-05-09 -> 2024-05-10
$ws.Range("G2:G255").Value = 45422

# Row-level balance corrections (r -> Saldo Previsto, Vl. Projetado, Vl. Total)
$rowUpdates = @(
    @{ Row = 2;   D = 122238.06; E = 0;       H = 122238.06 },
    @{ Row = 5;   D = 4678.16;   E = 0;       H = 4678.16 },
    @{ Row = 8;   D = 6012.34;   E = 0;       H = 6012.34 },
    @{ Row = 15;  D = 17807.2;   E = 0;       H = 17807.2 },
    @{ Row = 17;  D = 5976.37;   E = 0;       H = 5976.37 },
    @{ Row = 42;  D = 7609.73;   E = 0;       H = 7609.73 },
    @{ Row = 48;  D = 932.16;    E = 0;       H = 932.16 },
    @{ Row = 57;  D = 2493.1;    E = 1043.08; H = 3536.18 },
    @{ Row = 59;  D = 7194.4;    E = 0;       H = 7194.4 },
    @{ Row = 98;  D = 8380.14;   E = 0;       H = 8380.14 },
    @{ Row = 103; D = 22792.06;  E = 0;       H = 22792.06 },
    @{ Row = 107; D = 26843.76;  E = 0;       H = 26843.76 },
    @{ Row = 131; D = 4219.04;   E = 0;       H = 4219.04 },
    @{ Row = 141; D = 81519.55;  E = 0;       H = 81519.55 },
    @{ Row = 155; D = 1264.14;   E = 0;       H = 1264.14 },
    @{ Row = 167; D = 1782.11;   E = 0;       H = 1782.11 },
    @{ Row = 220; D = 25413.66;  E = 0;       H = 25413.66 },
    @{ Row = 224; D = 6742.4;    E = 0;       H = 6742.4 },
    @{ Row = 238; D = 8452.2;    E = 0;       H = 8452.2 },
    @{ Row = 243; D = 381.93;    E = 0;       H = 381.93 }
)

foreach ($u in $rowUpdates) {
    $r = $u.Row
    $ws.Range("D$r").Value = $u.D
    $ws.Range("E$r").Value = $u.E
    $ws.Range("H$r").Value = $u.H
}

# Move the active selection
$ws.Range("K15").Select() | Out-Null
